$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.212.02"
$ws.Range("D3").Value = "1.569.95"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.34%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "210.84"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +1.80%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.491"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +0.22%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "22.03"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("E9").Value = "  +0.22%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0598"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -0.04%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0869"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "1.793.07"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "1.569.31"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("E14").Value = "  +0.86%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.519"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "27.151.02"
$ws.Range("E16").Value = "  +0.66%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "62.24"
$cell.ClearFormats()
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "7.52"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +2.29%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "216.65"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "0.0₃0702"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("E22").Value = "  +0.99%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "9.22"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  +0.19%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "153.70"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +0.20%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "6.63"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  +2.42%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.0472"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "1.449.84"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("E35").Value = "  +7.10%  "
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  +2.36%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.810"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  +0.27%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.35"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("E44").Value = "  -0.78%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "64.42"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").Value = "1.704.79"
$ws.Range("E47").Value = "  +0.43%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "85.99"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("E50").Value = "  +0.43%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0956"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +0.00%  "

Write-Host "Done applying crypto updates"
